# Regen save_data: update column G ("K" = strikeouts) values for each row,
# calculated/re-derived from the source data (replacing the old "Strike#"-derived values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    3  = 0
    4  = 1
    5  = 1
    7  = 1
    8  = 4
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 3
    14 = 1
    15 = 2
    16 = 2
    17 = 3
    18 = 1
    19 = 0
    20 = 3
    21 = 0
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 2
    33 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
